$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "SearchTest" (1st sheet): move the selection to D9
# ---------------------------------------------------------------------
$search = $wb.Worksheets.Item("SearchTest")
$search.Range("D9").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "VillaTest" (2nd sheet): widen column B to fit its contents
# ---------------------------------------------------------------------
$villa = $wb.Worksheets.Item("VillaTest")
$villa.Columns.Item(2).AutoFit() | Out-Null

# ---------------------------------------------------------------------
# New sheet "ProductSearch" appended after VillaTest
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$productSearch = $wb.Worksheets.Add($null, $lastSheet)
$productSearch.Name = "ProductSearch"

$productSearch.Range("A1").Value = "Product Name"
$productSearch.Range("A2").Value = "OSSOBERRY Micro USB, USB Type C, Lightning OTG Adapter"
$productSearch.Range("A3").Value = "OSSOBERRY Micro USB, USB Type C, Lightning OTG Adapter (Pack of 1)"

$productSearch.Columns.Item(1).AutoFit() | Out-Null
$productSearch.Range("F17").Select() | Out-Null

# ---------------------------------------------------------------------
# New sheet "CardDetails" appended after ProductSearch
# ---------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$cardDetails = $wb.Worksheets.Add($null, $lastSheet2)
$cardDetails.Name = "CardDetails"

$cardDetails.Range("A1").Value = "Card"
$cardDetails.Range("B1").Value = "Card No"
$cardDetails.Range("C1").Value = "Expiry Date"
$cardDetails.Range("D1").Value = "CVV"

$cardDetails.Range("A2").Value = "Credit Card"

# Card No / Expiry Date are stored as Text so Excel doesn't reinterpret
# the long digit string or the "02/30" pattern. Expiry Date is written
# first so the shared-string table ends up in the same order the source
# workbook used.
$cardDetails.Range("C2").NumberFormat = "@"
$cardDetails.Range("C2").Value = "02/30"

$cardDetails.Range("B2").NumberFormat = "@"
$cardDetails.Range("B2").Value = "4315813955699002"

$cardDetails.Range("D2").Value = 143

$cardDetails.Columns.Item(2).AutoFit() | Out-Null
$cardDetails.Columns.Item(3).AutoFit() | Out-Null

$cardDetails.PageSetup.Orientation = 1

$cardDetails.Range("B2").Select() | Out-Null
